{"js": "// Change 1: \"\u03a4\u03b7\u03bd \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002\" -> \"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002\"\nconst body = context.document.body;\n\nconst search1 = body.search(\"\u03a4\u03b7\u03bd \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002\", { matchCase: true });\nsearch1.load(\"items\");\nawait context.sync();\n\nif (search1.items.length > 0) {\n  search1.items[0].insertText(\"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002\", \"Replace\");\n  await context.sync();\n}\n\n// Change 2: paragraph mentioning \u03a6.350.2/1/32958 gets replaced with the new decision text\n// (it spans 3 runs in the original, including a red-colored plain space run; the\n// whole paragraph text is replaced with a single plain run in the new content).\nconst search2 = body.search(\n  \"\u03a4\u03b7\u03bd \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.350.2/1/32958/\u03953/27-2-2018  (\u0391\u0394\u0391:6\u03a0414653\u03a0\u03a3-7\u0395\u039d) \u03a5\u03c0\u03bf\u03c5\u03c1\u03b3\u03b9\u03ba\u03ae \u0391\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7 \u03bc\u03b5 \u03b8\u03ad\u03bc\u03b1: \u00ab\u03a4\u03bf\u03c0\u03bf\u03b8\u03ad\u03c4\u03b7\u03c3\u03b7\",\n  { matchCase: true }\n);\nsearch2.load(\"items\");\nawait context.sync();\n\nif (search2.items.length > 0) {\n  const hitRange = search2.items[0];\n  const paragraph = hitRange.paragraphs.getFirst();\n  const paraRange = paragraph.getRange();\n  paraRange.insertText(\n    \"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.351.1/11/48020/\u03953/28-3-2019 (\u0391\u0394\u0391: \u03a9\u03a9\u03a4\u03974653\u03a0\u03a3-\u0392\u03943) \u03a5\u03c0\u03bf\u03c5\u03c1\u03b3\u03b9\u03ba\u03ae \u0391\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7 \u03bc\u03b5 \u03b8\u03ad\u03bc\u03b1: \u00ab\u03a4\u03bf\u03c0\u03bf\u03b8\u03ad\u03c4\u03b7\u03c3\u03b7 \u03a0\u03b5\u03c1\u03b9\u03c6\u03b5\u03c1\u03b5\u03b9\u03b1\u03ba\u03ce\u03bd \u0394\u03b9\u03b5\u03c5\u03b8\u03c5\u03bd\u03c4\u03ce\u03bd \u0395\u03ba\u03c0\u03b1\u03af\u03b4\u03b5\u03c5\u03c3\u03b7\u03c2\u00bb\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: \"\u03a4\u03b7\u03bd \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002\" -> \"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002\"\n$oldText1 = \"\u03a4\u03b7\u03bd \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002\"\n$newText1 = \"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.353.1/324/105657/\u03941/8-10-2002\"\n$find1 = $d.Content.Find\n$find1.Execute($oldText1, $false, $false, $false, $false, $false, $true, 1, $false, $newText1, 2)\n\n# Change 2: paragraph mentioning \u03a6.350.2/1/32958 gets replaced with the new decision text\n# (it spans 3 runs in the original, including a red-colored plain space run; the\n# whole paragraph text is replaced with a single plain run in the new content).\n$oldText2 = \"\u03a4\u03b7\u03bd \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.350.2/1/32958/\u03953/27-2-2018  (\u0391\u0394\u0391:6\u03a0414653\u03a0\u03a3-7\u0395\u039d) \u03a5\u03c0\u03bf\u03c5\u03c1\u03b3\u03b9\u03ba\u03ae \u0391\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7 \u03bc\u03b5 \u03b8\u03ad\u03bc\u03b1: \u00ab\u03a4\u03bf\u03c0\u03bf\u03b8\u03ad\u03c4\u03b7\u03c3\u03b7 \u03a0\u03b5\u03c1\u03b9\u03c6\u03b5\u03c1\u03b5\u03b9\u03b1\u03ba\u03ce\u03bd \u0394\u03b9\u03b5\u03c5\u03b8\u03c5\u03bd\u03c4\u03ce\u03bd \u0395\u03ba\u03c0\u03b1\u03af\u03b4\u03b5\u03c5\u03c3\u03b7\u03c2 \u00bb\"\n$newText2 = \"\u03a4\u03b7 \u03bc\u03b5 \u03b1\u03c1\u03b9\u03b8. \u03a6.351.1/11/48020/\u03953/28-3-2019 (\u0391\u0394\u0391: \u03a9\u03a9\u03a4\u03974653\u03a0\u03a3-\u0392\u03943) \u03a5\u03c0\u03bf\u03c5\u03c1\u03b3\u03b9\u03ba\u03ae \u0391\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7 \u03bc\u03b5 \u03b8\u03ad\u03bc\u03b1: \u00ab\u03a4\u03bf\u03c0\u03bf\u03b8\u03ad\u03c4\u03b7\u03c3\u03b7 \u03a0\u03b5\u03c1\u03b9\u03c6\u03b5\u03c1\u03b5\u03b9\u03b1\u03ba\u03ce\u03bd \u0394\u03b9\u03b5\u03c5\u03b8\u03c5\u03bd\u03c4\u03ce\u03bd \u0395\u03ba\u03c0\u03b1\u03af\u03b4\u03b5\u03c5\u03c3\u03b7\u03c2\u00bb\"\n$find2 = $d.Content.Find\n$find2.Execute($oldText2, $false, $false, $false, $false, $false, $true, 1, $false, $newText2, 2)\n"}
